$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (column D) / 1h-volume (column E) figures, taken from the
# refreshed cryptos feed. Each row keeps the exact text the source snapshot
# uses (including the 2-space padding around the % in column E).
$updates = @(
    @{ Cell = 'D2'; Text = '34.135.06' },
    @{ Cell = 'E2'; Text = '  +0.01%  ' },
    @{ Cell = 'D3'; Text = '1.781.97' },
    @{ Cell = 'E3'; Text = '  -0.51%  ' },
    @{ Cell = 'D5'; Text = '225.81' },
    @{ Cell = 'E5'; Text = '  -0.57%  ' },
    @{ Cell = 'E6'; Text = '  +0.00%  ' },
    @{ Cell = 'E7'; Text = '  +0.28%  ' },
    @{ Cell = 'D8'; Text = '32.11' },
    @{ Cell = 'E8'; Text = '  -0.97%  ' },
    @{ Cell = 'E9'; Text = '  -1.03%  ' },
    @{ Cell = 'E10'; Text = '  -0.49%  ' },
    @{ Cell = 'D11'; Text = '0.0950' },
    @{ Cell = 'E11'; Text = '  +0.90%  ' },
    @{ Cell = 'D12'; Text = '2.038.19' },
    @{ Cell = 'E12'; Text = '  -0.53%  ' },
    @{ Cell = 'D13'; Text = '1.784.46' },
    @{ Cell = 'E13'; Text = '  -0.53%  ' },
    @{ Cell = 'D14'; Text = '10.92' },
    @{ Cell = 'E14'; Text = '  -4.95%  ' },
    @{ Cell = 'E15'; Text = '  -0.16%  ' },
    @{ Cell = 'D16'; Text = '34.110.16' },
    @{ Cell = 'E16'; Text = '  -0.01%  ' },
    @{ Cell = 'E17'; Text = '  -0.11%  ' },
    @{ Cell = 'D18'; Text = '67.57' },
    @{ Cell = 'E18'; Text = '  -0.70%  ' },
    @{ Cell = 'D19'; Text = '245.35' },
    @{ Cell = 'E19'; Text = '  +0.47%  ' },
    @{ Cell = 'D20'; Text = '0.0₃0787' },
    @{ Cell = 'E20'; Text = '  +1.05%  ' },
    @{ Cell = 'E21'; Text = '  +0.34%  ' },
    @{ Cell = 'D22'; Text = '10.87' },
    @{ Cell = 'E22'; Text = '  -0.54%  ' },
    @{ Cell = 'E23'; Text = '  -0.12%  ' },
    @{ Cell = 'E24'; Text = '  -1.06%  ' },
    @{ Cell = 'D25'; Text = '162.24' },
    @{ Cell = 'E25'; Text = '  +0.35%  ' },
    @{ Cell = 'E26'; Text = '  -0.76%  ' },
    @{ Cell = 'D27'; Text = '16.28' },
    @{ Cell = 'E27'; Text = '  +0.07%  ' },
    @{ Cell = 'E28'; Text = '  +0.15%  ' },
    @{ Cell = 'E29'; Text = '  +0.43%  ' },
    @{ Cell = 'E30'; Text = '  -0.87%  ' },
    @{ Cell = 'E31'; Text = '  -0.47%  ' },
    @{ Cell = 'E32'; Text = '  +1.13%  ' },
    @{ Cell = 'E33'; Text = '  +2.59%  ' },
    @{ Cell = 'E34'; Text = '  -2.25%  ' },
    @{ Cell = 'D35'; Text = '1.447.43' },
    @{ Cell = 'E35'; Text = '  +2.90%  ' },
    @{ Cell = 'E36'; Text = '  +5.82%  ' },
    @{ Cell = 'D37'; Text = '0.651' },
    @{ Cell = 'E37'; Text = '  -0.82%  ' },
    @{ Cell = 'E38'; Text = '  +0.77%  ' },
    @{ Cell = 'E39'; Text = '  -0.58%  ' },
    @{ Cell = 'D40'; Text = '81.41' },
    @{ Cell = 'E40'; Text = '  +1.49%  ' },
    @{ Cell = 'E41'; Text = '  +1.41%  ' },
    @{ Cell = 'E42'; Text = '  +0.07%  ' },
    @{ Cell = 'D43'; Text = '0.912' },
    @{ Cell = 'E43'; Text = '  -1.20%  ' },
    @{ Cell = 'D44'; Text = '13.59' },
    @{ Cell = 'E44'; Text = '  +1.56%  ' },
    @{ Cell = 'E45'; Text = '  +2.80%  ' },
    @{ Cell = 'E46'; Text = '  -0.47%  ' },
    @{ Cell = 'E47'; Text = '  +0.33%  ' },
    @{ Cell = 'D48'; Text = '1.938.20' },
    @{ Cell = 'E48'; Text = '  -0.63%  ' },
    @{ Cell = 'E49'; Text = '  -6.71%  ' },
    @{ Cell = 'D50'; Text = '104.64' },
    @{ Cell = 'E50'; Text = '  -2.67%  ' },
    @{ Cell = 'E51'; Text = '  +0.33%  ' }
)

foreach ($u in $updates) {
    $value = $u.Text
    # Column D holds price text such as "225.81" or "34.135.06". Excel would
    # otherwise reinterpret a value like "225.81" as a Number, so prefix it with
    # a leading apostrophe (exactly what typing '225.81 into a cell does) to
    # force it to stay plain text, matching the source workbook.
    if ($u.Cell -like 'D*') {
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}
